# Auto-generated edit script applying the Cerberus_Profits.xlsx numeric updates
# across all 8 leveling sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ALC ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("H40").Value = 4998.25  # was 4997.5
$ws.Range("I40").Value = 4998.25  # was 4997.5
$ws.Range("K40").Value = 4998.25  # was 4997.5
$ws.Range("M40").Value = -4823.25  # was -4822.5
$ws.Range("H97").Value = 2889  # was 3328.6667
$ws.Range("J97").Value = 3016  # was 3594.4
$ws.Range("L97").Value = 9048  # was 10783.2
$ws.Range("N97").Value = -10040  # was -11775.2
$ws.Range("H106").Value = 9594.75  # was 10214.467
$ws.Range("I106").Value = 9867.799999999999  # was 10551.286
$ws.Range("K106").Value = 9867.799999999999  # was 10551.286
$ws.Range("M106").Value = -9236.799999999999  # was -9920.286
$ws.Range("H135").Value = 2188.8572  # was 1990.1875
$ws.Range("I135").Value = 1422.2727  # was 1353.75
$ws.Range("J135").Value = 4999.6665  # was 3899.5
$ws.Range("K135").Value = 12800.4543  # was 12183.75
$ws.Range("L135").Value = 44996.9985  # was 35095.5
$ws.Range("M135").Value = -10265.4543  # was -9648.75
$ws.Range("N135").Value = -50066.9985  # was -40165.5

# ---- Sheet 2: ARM ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1490.8889  # was 1560.5
$ws.Range("I2").Value = 1415.3572  # was 1495.5834
$ws.Range("K2").Value = 1415.3572  # was 1495.5834
$ws.Range("M2").Value = -1302.3572  # was -1382.5834
$ws.Range("H43").Value = 36998.5  # was 36999
$ws.Range("I43").Value = 36997  # was 0
$ws.Range("J43").Value = 37000  # was 36999
$ws.Range("K43").Value = 36997  # was 0
$ws.Range("L43").Value = 37000  # was 36999
$ws.Range("M43").Value = -36684  # was empty
$ws.Range("N43").Value = -37626  # was -37625
$ws.Range("H61").Value = 5447.75  # was 5315
$ws.Range("I61").Value = 3175.5652  # was 3109.8333
$ws.Range("K61").Value = 3175.5652  # was 3109.8333
$ws.Range("M61").Value = -2963.5652  # was -2897.8333
$ws.Range("H63").Value = 1703.1428  # was 1868.4
$ws.Range("I63").Value = 1487  # was 1585.5
$ws.Range("K63").Value = 1487  # was 1585.5
$ws.Range("M63").Value = -801  # was -899.5
$ws.Range("H66").Value = 1703.1428  # was 1868.4
$ws.Range("I66").Value = 1487  # was 1585.5
$ws.Range("K66").Value = 7435  # was 7927.5
$ws.Range("M66").Value = -4003  # was -4495.5
$ws.Range("H116").Value = 1490.8889  # was 1560.5
$ws.Range("I116").Value = 1415.3572  # was 1495.5834
$ws.Range("K116").Value = 1415.3572  # was 1495.5834
$ws.Range("M116").Value = 878.6428000000001  # was 798.4166
$ws.Range("H135").Value = 0  # was 70000
$ws.Range("J135").Value = 0  # was 70000
$ws.Range("L135").Value = 0  # was 70000
$ws.Range("N135").ClearContents()  # was -80140
$ws.Range("H136").Value = 5447.75  # was 5315
$ws.Range("I136").Value = 3175.5652  # was 3109.8333
$ws.Range("K136").Value = 9526.695599999999  # was 9329.499899999999
$ws.Range("M136").Value = -6976.695599999999  # was -6779.499899999999

# ---- Sheet 3: BSM ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1490.8889  # was 1560.5
$ws.Range("I3").Value = 1415.3572  # was 1495.5834
$ws.Range("K3").Value = 1415.3572  # was 1495.5834
$ws.Range("M3").Value = -1301.3572  # was -1381.5834
$ws.Range("H82").Value = 42479.25  # was 42480.5
$ws.Range("J82").Value = 74958.5  # was 74961
$ws.Range("L82").Value = 74958.5  # was 74961
$ws.Range("N82").Value = -75724.5  # was -75727
$ws.Range("H85").Value = 42479.25  # was 42480.5
$ws.Range("J85").Value = 74958.5  # was 74961
$ws.Range("L85").Value = 74958.5  # was 74961
$ws.Range("N85").Value = -77610.5  # was -77613
$ws.Range("H86").Value = 10500.929  # was 7786.3335
$ws.Range("I86").Value = 3032.111  # was 2736.8125
$ws.Range("K86").Value = 3032.111  # was 2736.8125
$ws.Range("M86").Value = -1909.111  # was -1613.8125
$ws.Range("H89").Value = 10500.929  # was 7786.3335
$ws.Range("I89").Value = 3032.111  # was 2736.8125
$ws.Range("K89").Value = 15160.555  # was 13684.0625
$ws.Range("M89").Value = -9544.555  # was -8068.0625
$ws.Range("H94").Value = 10334.952  # was 10626.7
$ws.Range("I94").Value = 649.1177  # was 408.4375
$ws.Range("K94").Value = 649.1177  # was 408.4375
$ws.Range("M94").Value = -198.1177  # was 42.5625
$ws.Range("H99").Value = 2172.5  # was 0
$ws.Range("I99").Value = 2345  # was 0
$ws.Range("J99").Value = 2000  # was 0
$ws.Range("K99").Value = 2345  # was 0
$ws.Range("L99").Value = 2000  # was 0
$ws.Range("M99").Value = -847  # was empty
$ws.Range("N99").Value = -4996  # was empty

# ---- Sheet 4: CRP ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("H58").Value = 2945.5557  # was 3001.5
$ws.Range("J58").Value = 6749  # was 11000
$ws.Range("L58").Value = 6749  # was 11000
$ws.Range("N58").Value = -7155  # was -11406
$ws.Range("H99").Value = 2583.4375  # was 2572.5881
$ws.Range("J99").Value = 2400  # was 2399.5
$ws.Range("L99").Value = 2400  # was 2399.5
$ws.Range("N99").Value = -5396  # was -5395.5
$ws.Range("H122").Value = 2489.2307  # was 2632.7144
$ws.Range("J122").Value = 2504.6667  # was 3003
$ws.Range("L122").Value = 7514.000100000001  # was 9009
$ws.Range("N122").Value = -12414.0001  # was -13909
$ws.Range("H126").Value = 2583.4375  # was 2572.5881
$ws.Range("J126").Value = 2400  # was 2399.5
$ws.Range("L126").Value = 7200  # was 7198.5
$ws.Range("N126").Value = -12140  # was -12138.5
$ws.Range("H132").Value = 5144.923  # was 3686.7407
$ws.Range("I132").Value = 5144.923  # was 3686.7407
$ws.Range("K132").Value = 15434.769  # was 11060.2221
$ws.Range("M132").Value = -12904.769  # was -8530.222099999999
$ws.Range("H134").Value = 5866.6855  # was 5764.8335
$ws.Range("I134").Value = 6111.3  # was 5985.129
$ws.Range("K134").Value = 18333.9  # was 17955.387
$ws.Range("M134").Value = -15798.9  # was -15420.387
$ws.Range("H136").Value = 2945.5557  # was 3001.5
$ws.Range("J136").Value = 6749  # was 11000
$ws.Range("L136").Value = 20247  # was 33000
$ws.Range("N136").Value = -25347  # was -38100

# ---- Sheet 5: CUL ----
$ws = $wb.Worksheets.Item(5)
$ws.Range("H44").Value = 2722.111  # was 1671.0714
$ws.Range("I44").Value = 750  # was 199.6
$ws.Range("J44").Value = 4299.8  # was 5349.75
$ws.Range("K44").Value = 2250  # was 598.8
$ws.Range("L44").Value = 12899.4  # was 16049.25
$ws.Range("M44").Value = -1852  # was -200.8
$ws.Range("N44").Value = -13695.4  # was -16845.25
$ws.Range("H131").Value = 3826.1282  # was 3822.3206
$ws.Range("I131").Value = 5000  # was 3250
$ws.Range("J131").Value = 3810.883  # was 3837.3816
$ws.Range("K131").Value = 15000  # was 9750
$ws.Range("L131").Value = 11432.649  # was 11512.1448
$ws.Range("M131").Value = -9960  # was -4710
$ws.Range("N131").Value = -21512.649  # was -21592.1448
$ws.Range("H132").Value = 7298.4375  # was 7219.357
$ws.Range("J132").Value = 7856.636  # was 7857.6665
$ws.Range("L132").Value = 70709.724  # was 70718.9985
$ws.Range("N132").Value = -75769.724  # was -75778.9985

# ---- Sheet 6: GSM ----
$ws = $wb.Worksheets.Item(6)
$ws.Range("H24").Value = 111118.92  # was 169445.75
$ws.Range("J24").Value = 111118.92  # was 169445.75
$ws.Range("L24").Value = 111118.92  # was 169445.75
$ws.Range("N24").Value = -111464.92  # was -169791.75
$ws.Range("H97").Value = 2736.75  # was 3624.3125
$ws.Range("I97").Value = 1566.1818  # was 2071.9
$ws.Range("J97").Value = 5312  # was 6211.6665
$ws.Range("K97").Value = 1566.1818  # was 2071.9
$ws.Range("L97").Value = 5312  # was 6211.6665
$ws.Range("M97").Value = -1070.1818  # was -1575.9
$ws.Range("N97").Value = -6304  # was -7203.6665
$ws.Range("H102").Value = 2791.2  # was 2809.4092
$ws.Range("I102").Value = 2584.4211  # was 2600.4866
$ws.Range("K102").Value = 2584.4211  # was 2600.4866
$ws.Range("M102").Value = -962.4211  # was -978.4866000000002
$ws.Range("H126").Value = 6268.1177  # was 6268.2354
$ws.Range("I126").Value = 3766.6667  # was 3766.889
$ws.Range("K126").Value = 11300.0001  # was 11300.667
$ws.Range("M126").Value = -8830.000100000001  # was -8830.667000000001
$ws.Range("H132").Value = 3552.7144  # was 3126.2144
$ws.Range("I132").Value = 3505.6316  # was 3526.3157
$ws.Range("J132").Value = 4000  # was 2281.5557
$ws.Range("K132").Value = 10516.8948  # was 10578.9471
$ws.Range("L132").Value = 12000  # was 6844.6671
$ws.Range("M132").Value = -7986.8948  # was -8048.947100000001
$ws.Range("N132").Value = -17060  # was -11904.6671

# ---- Sheet 7: LTW ----
$ws = $wb.Worksheets.Item(7)
$ws.Range("H4").Value = 0  # was 1900
$ws.Range("J4").Value = 0  # was 1900
$ws.Range("L4").Value = 0  # was 1900
$ws.Range("N4").ClearContents()  # was -2126
$ws.Range("H5").Value = 3250  # was 1500000
$ws.Range("J5").Value = 3250  # was 1500000
$ws.Range("L5").Value = 3250  # was 1500000
$ws.Range("N5").Value = -3476  # was -1500226
$ws.Range("H23").Value = 12250  # was 10000
$ws.Range("I23").Value = 12250  # was 10000
$ws.Range("K23").Value = 12250  # was 10000
$ws.Range("M23").Value = -12020  # was -9770
$ws.Range("H28").Value = 0  # was 1900
$ws.Range("J28").Value = 0  # was 1900
$ws.Range("L28").Value = 0  # was 1900
$ws.Range("N28").ClearContents()  # was -2364
$ws.Range("H37").Value = 0  # was 1900
$ws.Range("J37").Value = 0  # was 1900
$ws.Range("L37").Value = 0  # was 1900
$ws.Range("N37").ClearContents()  # was -2114
$ws.Range("H108").Value = 41497.6  # was 44372.25
$ws.Range("J108").Value = 41497.6  # was 44372.25
$ws.Range("L108").Value = 41497.6  # was 44372.25
$ws.Range("N108").Value = -49177.6  # was -52052.25
$ws.Range("H122").Value = 3998.6667  # was 3999.6667
$ws.Range("I122").Value = 3998.4  # was 3999.5
$ws.Range("K122").Value = 11995.2  # was 11998.5
$ws.Range("M122").Value = -9545.200000000001  # was -9548.5
$ws.Range("H130").Value = 90332.336  # was 94999
$ws.Range("J130").Value = 90332.336  # was 94999
$ws.Range("L130").Value = 90332.336  # was 94999
$ws.Range("N130").Value = -100372.336  # was -105039
$ws.Range("H132").Value = 3367.5  # was 3225.1667
$ws.Range("I132").Value = 2597.25  # was 2607.1667
$ws.Range("J132").Value = 4291.8  # was 3843.1667
$ws.Range("K132").Value = 7791.75  # was 7821.500100000001
$ws.Range("L132").Value = 12875.4  # was 11529.5001
$ws.Range("M132").Value = -5261.75  # was -5291.500100000001
$ws.Range("N132").Value = -17935.4  # was -16589.5001
$ws.Range("H136").Value = 2583.2188  # was 2678.8
$ws.Range("I136").Value = 1497.8462  # was 1561.1818
$ws.Range("K136").Value = 4493.5386  # was 4683.5454
$ws.Range("M136").Value = -1943.5386  # was -2133.5454

# ---- Sheet 8: WVR ----
$ws = $wb.Worksheets.Item(8)
$ws.Range("H54").Value = 34998.375  # was 36664.11
$ws.Range("J54").Value = 49987  # was 49988.5
$ws.Range("L54").Value = 49987  # was 49988.5
$ws.Range("N54").Value = -51027  # was -51028.5
$ws.Range("H124").Value = 207714  # was 212874.5
$ws.Range("J124").Value = 207714  # was 212874.5
$ws.Range("L124").Value = 207714  # was 212874.5
$ws.Range("N124").Value = -217534  # was -222694.5
$ws.Range("H132").Value = 3437.1428  # was 3521.5
$ws.Range("I132").Value = 2843.3667  # was 2921.7932
$ws.Range("M132").Value = -6000.1001  # was -6235.3796
$ws.Range("H136").Value = 9315.852999999999  # was 7594.2095
$ws.Range("I136").Value = 7730.815  # was 6213.1714
$ws.Range("J136").Value = 15429.571  # was 13636.25
$ws.Range("K136").Value = 23192.445  # was 18639.5142
$ws.Range("L136").Value = 46288.713  # was 40908.75
$ws.Range("M136").Value = -20642.445  # was -16089.5142
$ws.Range("N136").Value = -51388.713  # was -46008.75
